$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row 4 (for the new taxon "C. bentleyi, involuta, and
#    striata"), pushing the existing Populus deltoides / blank / STACKS rows
#    down by one. Insert() copies the formatting of the row above, which
#    matches row 3's (Corallorhiza bentleyi) formatting -- exactly what the
#    new row needs.
# ---------------------------------------------------------------------------
$ws.Rows("4:4").Insert()

# ---------------------------------------------------------------------------
# 2. Row 4 (new taxon row: C. bentleyi, involuta, and striata). Populated
#    first so the new shared strings it introduces ("C. bentleyi, involuta,
#    and striata" and "2,721 (250 bp min)") get the earliest new string-table
#    slots, matching the order they were authored in.
# ---------------------------------------------------------------------------
$ws.Range("A4").Value2 = "C. bentleyi, involuta, and striata"
$ws.Range("B4").Value2 = "ISSRseq"
$ws.Range("C4").Value2 = "pooled PCRs, sheared"
$ws.Range("D4").Value2 = "de novo"
$ws.Range("E4").Value2 = "2,721 (250 bp min)"
$ws.Range("F4").Value2 = 1075811
$ws.Range("G4").Value2 = 12095
$ws.Range("H4").Value2 = 8542
$ws.Range("I4").Value2 = 492
$ws.Range("J4").Value2 = 0

# ---------------------------------------------------------------------------
# 3. Update the header row (row 1). Column H's label changes, and two new
#    columns (I, J) are introduced.
# ---------------------------------------------------------------------------
$ws.Range("J1").Value2 = "found in all individuals"
$ws.Range("I1").Value2 = "found in 90% of individuals"
$ws.Range("H1").Value2 = "filtered SNPs"

# new columns I & J need the same centered/general formatting as the rest
# of the header row
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4. Row 2 (Asarum canadense) -- add new "found in 90%"/"found in all" data.
# ---------------------------------------------------------------------------
$ws.Range("I2").Value2 = 138
$ws.Range("J2").Value2 = 4

# ---------------------------------------------------------------------------
# 5. Row 3 (Corallorhiza bentleyi) -- add new "found in 90%"/"found in all".
# ---------------------------------------------------------------------------
$ws.Range("I3").Value2 = 1225
$ws.Range("J3").Value2 = 174

# ---------------------------------------------------------------------------
# 6. Row 5 (Populus deltoides, shifted down from old row 4) -- add new
#    "found in 90%"/"found in all" columns as "NA" text (matching the
#    existing "NA" entries in columns E/F on that row).
# ---------------------------------------------------------------------------
$ws.Range("E5").Copy() | Out-Null
$ws.Range("I5:J5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("I5").Value2 = "NA"
$ws.Range("J5").Value2 = "NA"

# ---------------------------------------------------------------------------
# 7. Row 10 (ISSRseq / Cypripedium macranthum, shifted down from old row 9)
#    -- add new "found in 90%"/"found in all" data.
# ---------------------------------------------------------------------------
$ws.Range("I10").Value2 = 885
$ws.Range("J10").Value2 = 648

# ---------------------------------------------------------------------------
# 8. Column widths: column A widened, column I (existing) widened, and new
#    column J given an explicit width. (The engine quantizes ColumnWidth to
#    1/6-character increments on save, so the inputs below are chosen to
#    land as close as the quantizer allows to the saved widths of 30,
#    25.28515625 and 21.5703125 characters respectively.)
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 29.1
$ws.Columns("I").ColumnWidth = 24.5
$ws.Columns("J").ColumnWidth = 20.65

# ---------------------------------------------------------------------------
# 9. Selection / active cell, matching the saved view state in the diff.
# ---------------------------------------------------------------------------
$ws.Range("H4").Select()
